$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.949.44'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = '1.912.07'
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = "'324.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("D7").Value = "'0.4592"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("D8").Value = "'0.3820"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.96%  '
$ws.Range("D9").Value = "'0.07690"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.54%  '
$ws.Range("D10").Value = "'0.9797"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("E11").Value = '  -0.88%  '
$ws.Range("D12").Value = '1.905.33'
$ws.Range("E12").Value = '  -2.08%  '
$ws.Range("D13").Value = "'5.681"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.31%  '
$ws.Range("D14").Value = "'6.947"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.36%  '
$ws.Range("D15").Value = "'0.07037"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.72%  '
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("E17").Value = '  -2.85%  '
$ws.Range("D18").Value = "'0.000009421"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.51%  '
$ws.Range("E19").Value = '  -2.13%  '
$ws.Range("D20").Value = "'0.9998"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").Value = '28.940.28'
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").Value = "'5.315"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.80%  '
$ws.Range("E23").Value = '  -1.36%  '
$ws.Range("E24").Value = '  -0.83%  '
$ws.Range("D25").Value = "'158.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.60%  '
$ws.Range("E26").Value = '  -1.76%  '
$ws.Range("D27").Value = "'5.679"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.63%  '
$ws.Range("D28").Value = "'117.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.90%  '
$ws.Range("E29").Value = '  +2.53%  '
$ws.Range("D30").Value = "'0.09304"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("D31").Value = "'0.8642"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.50%  '
$ws.Range("E32").Value = '  -1.00%  '
$ws.Range("D33").Value = "'1.243"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.78%  '
$ws.Range("D34").Value = "'3.055"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.76%  '
$ws.Range("D35").Value = "'0.05702"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.08%  '
$ws.Range("D36").Value = "'1.157"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("D37").Value = "'0.9997"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.32%  '
$ws.Range("D38").Value = "'0.02039"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.60%  '
$ws.Range("D39").Value = "'7.486"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.35%  '
$ws.Range("D40").Value = "'0.5499"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.28%  '
$ws.Range("D41").Value = "'2.957"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.67%  '
$ws.Range("E42").Value = '  -1.26%  '
$ws.Range("D43").Value = "'9.372"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.65%  '
$ws.Range("D44").Value = "'0.000002846"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.26%  '
$ws.Range("D45").Value = "'2.177"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.97%  '
$ws.Range("D46").Value = "'0.5174"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.54%  '
$ws.Range("D47").Value = "'11.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.60%  '
$ws.Range("D48").Value = "'0.06886"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("D49").Value = "'1.780"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.16%  '
$ws.Range("D50").Value = "'110.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("E51").Value = '  -0.39%  '
